# Remove the guide entries that are "above" the earlier cluster and have
# no upstream match (rows with start positions 125424506 / 125424552 /
# 125424557 / 125424558). These four rows were duplicated in the sheet
# (rows 15-18 and again at rows 33-36); delete both occurrences.
#
# Delete from the bottom up so earlier row numbers stay valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33:A36").EntireRow.Delete()
$ws.Range("A15:A18").EntireRow.Delete()
